$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.32"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.74"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.275"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05730"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.434"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8092"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8750"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1424"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07376"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03018"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03126"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09388"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.936"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001586"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04801"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005839"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006142"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005098"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009978"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001500"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.749"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.299"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.188"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1320"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001100"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03904"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006783"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1067"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002530"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007476"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005639"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5999"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1742"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01010"
